# dados/BIBI/Dados_BIBI_PF/faturamento_diario.xlsx
# "terminei vendas BIBI e arrumei vendas ADD na analise e no dash"
#
# April (C=4) was missing day 29: the last April row (day 28, row 29) had a
# placeholder/old value, and everything from "March day 1" onward (old row
# 30) needs to shift down by one row to make room for the newly completed
# April day-29 entry. Net effect on the sheet:
#   - B29 (April, day 28) total_venda corrected: 16598.71 -> 25680.18
#   - a new row is inserted at row 30 for April day 29 (value 21976.8),
#     pushing every subsequent row (old 30..118) down by one (new 31..119)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the April day-28 total_venda value.
$ws.Cells.Item(29, 2).Value = 25680.18

# Insert the new row for April day 29, shifting the rest of the table down.
$ws.Rows.Item(30).Insert()

$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = 21976.8
$ws.Cells.Item(30, 3).Value = 4
$ws.Cells.Item(30, 4).Value = 2025
$ws.Cells.Item(30, 5).Value = "04/2025"
